# Atualização de bases das ligas, do dia: 24-02-2024 às 12:40
# This script swaps/updates rows in the "Chile Primera Division" sheet
# to reflect refreshed match-odds data, per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows [62, 63] ---
$ws.Range("B62").Value = 6078868
$ws.Range("C62").Value = "Chile Primera Division"
$ws.Range("D62").Value = "Chile Primera Division"
$ws.Range("E62").Value = 45003.75
$ws.Range("F62").Value = "Cobresal"
$ws.Range("G62").Value = "Colo Colo"
$ws.Range("H62").Value = 3
$ws.Range("I62").Value = 1
$ws.Range("J62").Value = "H"
$ws.Range("K62").Value = 2.875
$ws.Range("L62").Value = 3.25
$ws.Range("M62").Value = 2.25
$ws.Range("N62").Value = 3.3
$ws.Range("O62").Value = 3.3
$ws.Range("P62").Value = 2.25
$ws.Range("Q62").Value = 0.25
$ws.Range("R62").Value = 1.875
$ws.Range("S62").Value = 1.925
$ws.Range("T62").Value = 2.5
$ws.Range("U62").Value = 2
$ws.Range("V62").Value = 1.8
$ws.Range("W62").Value = 2.3
$ws.Range("X62").Value = -1
$ws.Range("Y62").Value = -1
$ws.Range("Z62").Value = 0.875
$ws.Range("AA62").Value = -1
$ws.Range("AB62").Value = 1
$ws.Range("AC62").Value = -1
$ws.Range("B63").Value = 6078944
$ws.Range("C63").Value = "Chile Primera Division"
$ws.Range("D63").Value = "Chile Primera Division"
$ws.Range("E63").Value = 45003.75
$ws.Range("F63").Value = "OHiggins"
$ws.Range("G63").Value = "Coquimbo Unido"
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 2
$ws.Range("J63").Value = "A"
$ws.Range("K63").Value = 2
$ws.Range("L63").Value = 3.25
$ws.Range("M63").Value = 3.4
$ws.Range("N63").Value = 2.15
$ws.Range("O63").Value = 3.25
$ws.Range("P63").Value = 3.6
$ws.Range("Q63").Value = -0.25
$ws.Range("R63").Value = 1.8
$ws.Range("S63").Value = 2
$ws.Range("T63").Value = 2.25
$ws.Range("U63").Value = 1.85
$ws.Range("V63").Value = 1.95
$ws.Range("W63").Value = -1
$ws.Range("X63").Value = -1
$ws.Range("Y63").Value = 2.6
$ws.Range("Z63").Value = -1
$ws.Range("AA63").Value = 1
$ws.Range("AB63").Value = -0.5
$ws.Range("AC63").Value = 0.475

# --- Swap rows [179, 180] ---
$ws.Range("B179").Value = 7157967
$ws.Range("C179").Value = "Chile Primera Division"
$ws.Range("D179").Value = "Chile Primera Division"
$ws.Range("E179").Value = 45183.79166666666
$ws.Range("F179").Value = "Huachipato"
$ws.Range("G179").Value = "Palestino"
$ws.Range("H179").Value = 2
$ws.Range("I179").Value = 2
$ws.Range("J179").Value = "D"
$ws.Range("K179").Value = 2.375
$ws.Range("L179").Value = 3.2
$ws.Range("M179").Value = 3
$ws.Range("N179").Value = 2.75
$ws.Range("O179").Value = 3.2
$ws.Range("P179").Value = 2.7
$ws.Range("Q179").Value = 0
$ws.Range("R179").Value = 1.925
$ws.Range("S179").Value = 1.875
$ws.Range("T179").Value = 2.5
$ws.Range("U179").Value = 2
$ws.Range("V179").Value = 1.8
$ws.Range("W179").Value = -1
$ws.Range("X179").Value = 2.2
$ws.Range("Y179").Value = -1
$ws.Range("Z179").Value = 0
$ws.Range("AA179").Value = -0
$ws.Range("AB179").Value = 1
$ws.Range("AC179").Value = -1
$ws.Range("B180").Value = 7082624
$ws.Range("C180").Value = "Chile Primera Division"
$ws.Range("D180").Value = "Chile Primera Division"
$ws.Range("E180").Value = 45183.79166666666
$ws.Range("F180").Value = "Colo Colo"
$ws.Range("G180").Value = "Deportes Copiapo"
$ws.Range("H180").Value = 1
$ws.Range("I180").Value = 1
$ws.Range("J180").Value = "D"
$ws.Range("K180").Value = 1.333
$ws.Range("L180").Value = 5
$ws.Range("M180").Value = 8
$ws.Range("N180").Value = 1.45
$ws.Range("O180").Value = 4.75
$ws.Range("P180").Value = 7
$ws.Range("Q180").Value = -1.25
$ws.Range("R180").Value = 1.925
$ws.Range("S180").Value = 1.875
$ws.Range("T180").Value = 3
$ws.Range("U180").Value = 1.875
$ws.Range("V180").Value = 1.925
$ws.Range("W180").Value = -1
$ws.Range("X180").Value = 3.75
$ws.Range("Y180").Value = -1
$ws.Range("Z180").Value = -1
$ws.Range("AA180").Value = 0.875
$ws.Range("AB180").Value = -1
$ws.Range("AC180").Value = 0.925

# --- Swap rows [220, 221] ---
$ws.Range("B220").Value = 6077497
$ws.Range("C220").Value = "Chile Primera Division"
$ws.Range("D220").Value = "Chile Primera Division"
$ws.Range("E220").Value = 45256.85416666666
$ws.Range("F220").Value = "Deportes Copiapo"
$ws.Range("G220").Value = "Nublense"
$ws.Range("H220").Value = 1
$ws.Range("I220").Value = 1
$ws.Range("J220").Value = "D"
$ws.Range("K220").Value = 2.6
$ws.Range("L220").Value = 3.4
$ws.Range("M220").Value = 2.6
$ws.Range("N220").Value = 2.8
$ws.Range("O220").Value = 3.2
$ws.Range("P220").Value = 2.7
$ws.Range("Q220").Value = 0
$ws.Range("R220").Value = 1.95
$ws.Range("S220").Value = 1.9
$ws.Range("T220").Value = 2.25
$ws.Range("U220").Value = 2
$ws.Range("V220").Value = 1.85
$ws.Range("W220").Value = -1
$ws.Range("X220").Value = 2.2
$ws.Range("Y220").Value = -1
$ws.Range("Z220").Value = 0
$ws.Range("AA220").Value = -0
$ws.Range("AB220").Value = -0.5
$ws.Range("AC220").Value = 0.425
$ws.Range("B221").Value = 6077763
$ws.Range("C221").Value = "Chile Primera Division"
$ws.Range("D221").Value = "Chile Primera Division"
$ws.Range("E221").Value = 45256.85416666666
$ws.Range("F221").Value = "Curico Unido"
$ws.Range("G221").Value = "Magallanes"
$ws.Range("H221").Value = 3
$ws.Range("I221").Value = 4
$ws.Range("J221").Value = "A"
$ws.Range("K221").Value = 2.15
$ws.Range("L221").Value = 3.5
$ws.Range("M221").Value = 3.2
$ws.Range("N221").Value = 2.625
$ws.Range("O221").Value = 3.5
$ws.Range("P221").Value = 2.6
$ws.Range("Q221").Value = 0
$ws.Range("R221").Value = 1.95
$ws.Range("S221").Value = 1.9
$ws.Range("T221").Value = 2.75
$ws.Range("U221").Value = 1.975
$ws.Range("V221").Value = 1.875
$ws.Range("W221").Value = -1
$ws.Range("X221").Value = -1
$ws.Range("Y221").Value = 1.6
$ws.Range("Z221").Value = -1
$ws.Range("AA221").Value = 0.8999999999999999
$ws.Range("AB221").Value = 0.9750000000000001
$ws.Range("AC221").Value = -1

# --- Swap rows [224, 225] ---
$ws.Range("B224").Value = 6078265
$ws.Range("C224").Value = "Chile Primera Division"
$ws.Range("D224").Value = "Chile Primera Division"
$ws.Range("E224").Value = 45262.75
$ws.Range("F224").Value = "Audax Italiano"
$ws.Range("G224").Value = "Magallanes"
$ws.Range("H224").Value = 0
$ws.Range("I224").Value = 2
$ws.Range("J224").Value = "A"
$ws.Range("K224").Value = 1.666
$ws.Range("L224").Value = 3.75
$ws.Range("M224").Value = 5
$ws.Range("N224").Value = 2.25
$ws.Range("O224").Value = 3.3
$ws.Range("P224").Value = 3.3
$ws.Range("Q224").Value = -0.25
$ws.Range("R224").Value = 1.95
$ws.Range("S224").Value = 1.85
$ws.Range("T224").Value = 2.5
$ws.Range("U224").Value = 1.8
$ws.Range("V224").Value = 2
$ws.Range("W224").Value = -1
$ws.Range("X224").Value = -1
$ws.Range("Y224").Value = 2.3
$ws.Range("Z224").Value = -1
$ws.Range("AA224").Value = 0.8500000000000001
$ws.Range("AB224").Value = -1
$ws.Range("AC224").Value = 1
$ws.Range("B225").Value = 6077498
$ws.Range("C225").Value = "Chile Primera Division"
$ws.Range("D225").Value = "Chile Primera Division"
$ws.Range("E225").Value = 45262.75
$ws.Range("F225").Value = "Universidad Catolica"
$ws.Range("G225").Value = "Deportes Copiapo"
$ws.Range("H225").Value = 2
$ws.Range("I225").Value = 2
$ws.Range("J225").Value = "D"
$ws.Range("K225").Value = 1.65
$ws.Range("L225").Value = 3.8
$ws.Range("M225").Value = 5.25
$ws.Range("N225").Value = 1.909
$ws.Range("O225").Value = 3.6
$ws.Range("P225").Value = 4.2
$ws.Range("Q225").Value = -0.5
$ws.Range("R225").Value = 1.85
$ws.Range("S225").Value = 2
$ws.Range("T225").Value = 2.75
$ws.Range("U225").Value = 2.025
$ws.Range("V225").Value = 1.825
$ws.Range("W225").Value = -1
$ws.Range("X225").Value = 2.6
$ws.Range("Y225").Value = -1
$ws.Range("Z225").Value = -1
$ws.Range("AA225").Value = 1
$ws.Range("AB225").Value = 1.025
$ws.Range("AC225").Value = -1

# --- Swap rows [230, 231, 232] ---
$ws.Range("B230").Value = 6078267
$ws.Range("C230").Value = "Chile Primera Division"
$ws.Range("D230").Value = "Chile Primera Division"
$ws.Range("E230").Value = 45268.75
$ws.Range("F230").Value = "Huachipato"
$ws.Range("G230").Value = "Audax Italiano"
$ws.Range("H230").Value = 2
$ws.Range("I230").Value = 0
$ws.Range("J230").Value = "H"
$ws.Range("K230").Value = 1.5
$ws.Range("L230").Value = 4.333
$ws.Range("M230").Value = 6
$ws.Range("N230").Value = 1.444
$ws.Range("O230").Value = 4.75
$ws.Range("P230").Value = 7
$ws.Range("Q230").Value = -1.25
$ws.Range("R230").Value = 2.025
$ws.Range("S230").Value = 1.825
$ws.Range("T230").Value = 2.75
$ws.Range("U230").Value = 1.8
$ws.Range("V230").Value = 2.05
$ws.Range("W230").Value = 0.444
$ws.Range("X230").Value = -1
$ws.Range("Y230").Value = -1
$ws.Range("Z230").Value = 1.025
$ws.Range("AA230").Value = -1
$ws.Range("AB230").Value = -1
$ws.Range("AC230").Value = 1.05
$ws.Range("B231").Value = 6078997
$ws.Range("C231").Value = "Chile Primera Division"
$ws.Range("D231").Value = "Chile Primera Division"
$ws.Range("E231").Value = 45268.75
$ws.Range("F231").Value = "Union Espanola"
$ws.Range("G231").Value = "Cobresal"
$ws.Range("H231").Value = 1
$ws.Range("I231").Value = 0
$ws.Range("J231").Value = "H"
$ws.Range("K231").Value = 3.8
$ws.Range("L231").Value = 3.6
$ws.Range("M231").Value = 1.909
$ws.Range("N231").Value = 2.7
$ws.Range("O231").Value = 3.6
$ws.Range("P231").Value = 2.45
$ws.Range("Q231").Value = 0
$ws.Range("R231").Value = 1.975
$ws.Range("S231").Value = 1.825
$ws.Range("T231").Value = 2.75
$ws.Range("U231").Value = 1.775
$ws.Range("V231").Value = 2.025
$ws.Range("W231").Value = 1.7
$ws.Range("X231").Value = -1
$ws.Range("Y231").Value = -1
$ws.Range("Z231").Value = 0.9750000000000001
$ws.Range("AA231").Value = -1
$ws.Range("AB231").Value = -1
$ws.Range("AC231").Value = 1.025
$ws.Range("B232").Value = 6143704
$ws.Range("C232").Value = "Chile Primera Division"
$ws.Range("D232").Value = "Chile Primera Division"
$ws.Range("E232").Value = 45268.75
$ws.Range("F232").Value = "Curico Unido"
$ws.Range("G232").Value = "Colo Colo"
$ws.Range("H232").Value = 0
$ws.Range("I232").Value = 1
$ws.Range("J232").Value = "A"
$ws.Range("K232").Value = 6.5
$ws.Range("L232").Value = 4.75
$ws.Range("M232").Value = 1.4
$ws.Range("N232").Value = 12
$ws.Range("O232").Value = 8.5
$ws.Range("P232").Value = 1.166
$ws.Range("Q232").Value = 2
$ws.Range("R232").Value = 2
$ws.Range("S232").Value = 1.8
$ws.Range("T232").Value = 3.25
$ws.Range("U232").Value = 1.875
$ws.Range("V232").Value = 1.925
$ws.Range("W232").Value = -1
$ws.Range("X232").Value = -1
$ws.Range("Y232").Value = 0.1659999999999999
$ws.Range("Z232").Value = 1
$ws.Range("AA232").Value = -1
$ws.Range("AB232").Value = -1
$ws.Range("AC232").Value = 0.925

# --- Swap rows [233, 237] ---
$ws.Range("B233").Value = 6078269
$ws.Range("C233").Value = "Chile Primera Division"
$ws.Range("D233").Value = "Chile Primera Division"
$ws.Range("E233").Value = 45269.75
$ws.Range("F233").Value = "Universidad de Chile"
$ws.Range("G233").Value = "Nublense"
$ws.Range("H233").Value = 3
$ws.Range("I233").Value = 1
$ws.Range("J233").Value = "H"
$ws.Range("K233").Value = 1.85
$ws.Range("L233").Value = 3.4
$ws.Range("M233").Value = 4.333
$ws.Range("N233").Value = 1.8
$ws.Range("O233").Value = 3.6
$ws.Range("P233").Value = 4.5
$ws.Range("Q233").Value = -0.75
$ws.Range("R233").Value = 1.925
$ws.Range("S233").Value = 1.925
$ws.Range("T233").Value = 2.5
$ws.Range("U233").Value = 2.025
$ws.Range("V233").Value = 1.825
$ws.Range("W233").Value = 0.8
$ws.Range("X233").Value = -1
$ws.Range("Y233").Value = -1
$ws.Range("Z233").Value = 0.925
$ws.Range("AA233").Value = -1
$ws.Range("AB233").Value = 1.025
$ws.Range("AC233").Value = -1
$ws.Range("B237").Value = 6078268
$ws.Range("C237").Value = "Chile Primera Division"
$ws.Range("D237").Value = "Chile Primera Division"
$ws.Range("E237").Value = 45269.75
$ws.Range("F237").Value = "OHiggins"
$ws.Range("G237").Value = "Palestino"
$ws.Range("H237").Value = 0
$ws.Range("I237").Value = 1
$ws.Range("J237").Value = "A"
$ws.Range("K237").Value = 3.1
$ws.Range("L237").Value = 3.3
$ws.Range("M237").Value = 2.3
$ws.Range("N237").Value = 2.9
$ws.Range("O237").Value = 3.4
$ws.Range("P237").Value = 2.375
$ws.Range("Q237").Value = 0.25
$ws.Range("R237").Value = 1.8
$ws.Range("S237").Value = 2
$ws.Range("T237").Value = 2.75
$ws.Range("U237").Value = 2
$ws.Range("V237").Value = 1.8
$ws.Range("W237").Value = -1
$ws.Range("X237").Value = -1
$ws.Range("Y237").Value = 1.375
$ws.Range("Z237").Value = -1
$ws.Range("AA237").Value = 1
$ws.Range("AB237").Value = -1
$ws.Range("AC237").Value = 0.8

# --- Swap rows [246, 247] ---
$ws.Range("B246").Value = 7723522
$ws.Range("C246").Value = "Chile Primera Division"
$ws.Range("D246").Value = "Chile Primera Division"
$ws.Range("E246").Value = 45347.5
$ws.Range("F246").Value = "Huachipato"
$ws.Range("G246").Value = "Union La Calera"
$ws.Range("K246").Value = 2.1
$ws.Range("L246").Value = 3.4
$ws.Range("M246").Value = 3.5
$ws.Range("N246").Value = 1.85
$ws.Range("O246").Value = 3.6
$ws.Range("P246").Value = 4.2
$ws.Range("Q246").Value = -0.5
$ws.Range("R246").Value = 1.9
$ws.Range("S246").Value = 1.95
$ws.Range("T246").Value = 2.5
$ws.Range("U246").Value = 1.9
$ws.Range("V246").Value = 1.95
$ws.Range("W246").Value = 0
$ws.Range("X246").Value = 0
$ws.Range("Y246").Value = 0
$ws.Range("Z246").Value = 0
$ws.Range("AA246").Value = 0
$ws.Range("B247").Value = 7723523
$ws.Range("C247").Value = "Chile Primera Division"
$ws.Range("D247").Value = "Chile Primera Division"
$ws.Range("E247").Value = 45347.5
$ws.Range("F247").Value = "Deportes Iquique"
$ws.Range("G247").Value = "Everton de Vina"
$ws.Range("K247").Value = 3.3
$ws.Range("L247").Value = 3.3
$ws.Range("M247").Value = 2.2
$ws.Range("N247").Value = 2.9
$ws.Range("O247").Value = 3.3
$ws.Range("P247").Value = 2.45
$ws.Range("Q247").Value = 0
$ws.Range("R247").Value = 2.1
$ws.Range("S247").Value = 1.775
$ws.Range("T247").Value = 2.5
$ws.Range("U247").Value = 2
$ws.Range("V247").Value = 1.85
$ws.Range("W247").Value = 0
$ws.Range("X247").Value = 0
$ws.Range("Y247").Value = 0
$ws.Range("Z247").Value = 0
$ws.Range("AA247").Value = 0

# --- Standalone updates ---
$ws.Range("N244").Value = 1.727
$ws.Range("P244").Value = 4.75
$ws.Range("R244").Value = 1.95
$ws.Range("S244").Value = 1.9
$ws.Range("U244").Value = 1.9
$ws.Range("V244").Value = 1.95
$ws.Range("U245").Value = 2
$ws.Range("V245").Value = 1.85
$ws.Range("R248").Value = 2
$ws.Range("S248").Value = 2.25
